$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Training Schedule Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Training Schedule Overview")

$ws1.Range("A1").Value = "Product Development IMPLEMENTATION PROJECT - TRProductNING SCHEDULE"
$ws1.Range("C4").Value = "Enterprise Product Development Implementation"

$ws1.Range("A7").Value = "TRProductNING SCHEDULE SUMMARY"

$ws1.Range("A9").Value = "Product Development Fundamentals (Product-101)"
$ws1.Range("A10").Value = "Product Development Platform Overview (Product-102)"
$ws1.Range("A11").Value = "Data Analysis for Business (Product-201)"
$ws1.Range("A12").Value = "Advanced Product Techniques (Product-301)"
$ws1.Range("A13").Value = "ProductOps for IT Teams (Product-302)"
$ws1.Range("B13").Value = "Product Engineers, IT"
$ws1.Range("A14").Value = "Model Validation & QA (Product-303)"
$ws1.Range("B14").Value = "Product Engineers, QA"
$ws1.Range("A15").Value = "Executive Overview (Product-401)"
$ws1.Range("A16").Value = "Train-the-Trainer (Product-501)"

$ws1.Range("A18").Value = "TRProductNING SCHEDULE STATISTICS"

# Materialize the two newly-present (but empty) rows from the diff: row 6 and
# row 17. Touching OutlineLevel forces the row record to be written without
# adding height/visibility attributes, matching the target's bare <row r="n"/>.
$ws1.Rows.Item(6).OutlineLevel = 0
$ws1.Rows.Item(17).OutlineLevel = 0

# ---------------------------------------------------------------------------
# Sheet 2: "Detailed Training Schedule"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Detailed Training Schedule")

$ws2.Range("A1").Value = "DETProductLED TRProductNING SCHEDULE"

$ws2.Range("A4").Value = "Product-101"
$ws2.Range("B4").Value = "Product Development Fundamentals"
$ws2.Range("A5").Value = "Product-102"
$ws2.Range("B5").Value = "Product Development Platform Overview"
$ws2.Range("A6").Value = "Product-201"
$ws2.Range("A7").Value = "Product-201"
$ws2.Range("A8").Value = "Product-201"
$ws2.Range("A9").Value = "Product-301"
$ws2.Range("B9").Value = "Advanced Product Techniques"
$ws2.Range("A10").Value = "Product-301"
$ws2.Range("B10").Value = "Advanced Product Techniques"
$ws2.Range("A11").Value = "Product-302"
$ws2.Range("B11").Value = "ProductOps for IT Teams"
$ws2.Range("C11").Value = "Product Engineers, IT"
$ws2.Range("A12").Value = "Product-302"
$ws2.Range("B12").Value = "ProductOps for IT Teams"
$ws2.Range("C12").Value = "Product Engineers, IT"
$ws2.Range("A13").Value = "Product-303"
$ws2.Range("C13").Value = "Product Engineers, QA"
$ws2.Range("A14").Value = "Product-303"
$ws2.Range("C14").Value = "Product Engineers, QA"
$ws2.Range("A15").Value = "Product-401"
$ws2.Range("A16").Value = "Product-501"
$ws2.Range("A17").Value = "Product-501"
$ws2.Range("A18").Value = "Product-501"
$ws2.Range("A19").Value = "Product-501"
$ws2.Range("A20").Value = "Product-501"

$ws2.Rows.Item(2).OutlineLevel = 0

# ---------------------------------------------------------------------------
# Sheet 3: "Instructor Schedule"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Instructor Schedule")

$ws3.Range("B7").Value = "Advanced Product Techniques"
$ws3.Range("B8").Value = "Advanced Product Techniques"
$ws3.Range("B9").Value = "ProductOps for IT Teams"
$ws3.Range("B10").Value = "ProductOps for IT Teams"

$ws3.Rows.Item(2).OutlineLevel = 0

# ---------------------------------------------------------------------------
# Sheet 4: "Facility Schedule"
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Facility Schedule")

$ws4.Range("B7").Value = "Advanced Product Techniques"
$ws4.Range("B8").Value = "Advanced Product Techniques"
$ws4.Range("B9").Value = "ProductOps for IT Teams"
$ws4.Range("B10").Value = "ProductOps for IT Teams"

$ws4.Rows.Item(2).OutlineLevel = 0

# ---------------------------------------------------------------------------
# Sheet 5: "Participant Tracking"
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Participant Tracking")

$ws5.Range("F4").Value = "Product-101"
$ws5.Range("F5").Value = "Product-102"
$ws5.Range("F6").Value = "Product-301"
$ws5.Range("F7").Value = "Product-302"

$ws5.Range("E8").Value = "Product Engineer"
$ws5.Range("F8").Value = "Product-101"
$ws5.Range("E9").Value = "Product Engineer"
$ws5.Range("F9").Value = "Product-102"
$ws5.Range("E10").Value = "Product Engineer"
$ws5.Range("F10").Value = "Product-302"
$ws5.Range("E11").Value = "Product Engineer"
$ws5.Range("F11").Value = "Product-303"

$ws5.Range("F12").Value = "Product-101"
$ws5.Range("F13").Value = "Product-102"
$ws5.Range("F14").Value = "Product-401"

$ws5.Range("F15").Value = "Product-101"
$ws5.Range("F16").Value = "Product-102"
$ws5.Range("F17").Value = "Product-501"

$ws5.Range("F18").Value = "Product-101"
$ws5.Range("F19").Value = "Product-102"
$ws5.Range("F20").Value = "Product-301"
$ws5.Range("F21").Value = "Product-303"
$ws5.Range("F22").Value = "Product-501"

$ws5.Rows.Item(2).OutlineLevel = 0
